$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9394679069519043
$ws.Range("B1").Value = 2.822190046310425
$ws.Range("C1").Value = 5.406792163848877
$ws.Range("D1").Value = 2.092458009719849
$ws.Range("E1").Value = 1.181698322296143
